$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values that look numeric are prefixed with a literal
# leading apostrophe so Excel stores them as text (quotePrefix),
# matching the original text-formatted price strings exactly.

$ws.Range("D2").Value = "'63.090.54"
$ws.Range("E2").Value = "  -1.08%  "

$ws.Range("D3").Value = "'3.175.76"
$ws.Range("E3").Value = "  -4.46%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'592.13"
$ws.Range("E5").Value = "  -2.34%  "

$ws.Range("D6").Value = "'135.70"
$ws.Range("E6").Value = "  -4.72%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'3.172.38"
$ws.Range("E8").Value = "  -4.54%  "

$ws.Range("D9").Value = "'0.514"
$ws.Range("E9").Value = "  -1.17%  "

$ws.Range("E10").Value = "  -4.93%  "

$ws.Range("D11").Value = "'5.26"
$ws.Range("E11").Value = "  -4.78%  "

$ws.Range("D12").Value = "'0.455"
$ws.Range("E12").Value = "  -2.86%  "

$ws.Range("D13").Value = "'0.0000238"
$ws.Range("E13").Value = "  -3.81%  "

$ws.Range("D14").Value = "'34.74"
$ws.Range("E14").Value = "  -0.20%  "

$ws.Range("D15").Value = "'3.699.72"
$ws.Range("E15").Value = "  -4.50%  "

$ws.Range("D17").Value = "'3.178.41"
$ws.Range("E17").Value = "  -4.28%  "

$ws.Range("D18").Value = "'63.048.26"
$ws.Range("E18").Value = "  -1.23%  "

$ws.Range("D19").Value = "'6.60"
$ws.Range("E19").Value = "  -4.00%  "

$ws.Range("D20").Value = "'462.07"
$ws.Range("E20").Value = "  -3.85%  "

$ws.Range("D21").Value = "'13.93"
$ws.Range("E21").Value = "  -1.37%  "

$ws.Range("D22").Value = "'0.703"
$ws.Range("E22").Value = "  -4.69%  "

$ws.Range("D23").Value = "'7.64"
$ws.Range("E23").Value = "  -6.91%  "

$ws.Range("E24").Value = "  -2.25%  "

$ws.Range("D25").Value = "'83.63"
$ws.Range("E25").Value = "  -1.44%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").Value = "'2.68"
$ws.Range("E27").Value = "  -3.39%  "

$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").Value = "'7.76"
$ws.Range("E29").Value = "  -5.62%  "

$ws.Range("D30").Value = "'6.77"
$ws.Range("E30").Value = "  -6.46%  "

$ws.Range("D31").Value = "'2.03"
$ws.Range("E31").Value = "  -5.78%  "

$ws.Range("D32").Value = "'27.27"
$ws.Range("E32").Value = "  -5.92%  "

$ws.Range("E33").Value = "  -2.35%  "

$ws.Range("D34").Value = "'2.39"
$ws.Range("E34").Value = "  -5.79%  "

$ws.Range("E35").Value = "  -5.59%  "

$ws.Range("E36").Value = "  -3.77%  "

$ws.Range("D37").Value = "'51.34"
$ws.Range("E37").Value = "  -2.03%  "

$ws.Range("D38").Value = "'0.0₃0707"
$ws.Range("E38").Value = "  -5.19%  "

$ws.Range("D39").Value = "'0.0389"
$ws.Range("E39").Value = "  -3.11%  "

$ws.Range("D40").Value = "'405.91"
$ws.Range("E40").Value = "  -6.61%  "

$ws.Range("D41").Value = "'8.13"
$ws.Range("E41").Value = "  -2.54%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.113"
$ws.Range("E42").Value = "  -7.14%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.64"
$ws.Range("E43").Value = "  -4.87%  "

$ws.Range("D44").Value = "'2.798.59"
$ws.Range("E44").Value = "  -9.87%  "

$ws.Range("D45").Value = "'0.252"
$ws.Range("E45").Value = "  -4.29%  "

$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "'0.999"
$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "'2.13"
$ws.Range("E47").Value = "  -4.53%  "

$ws.Range("D48").Value = "'25.54"
$ws.Range("E48").Value = "  -3.13%  "

$ws.Range("D49").Value = "'124.21"
$ws.Range("E49").Value = "  -0.01%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.111"
$ws.Range("E50").Value = "  -2.28%  "

$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").Value = "'34.38"
$ws.Range("E51").Value = "  -7.81%  "

Write-Host "Done"